$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (header is row 1, data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

# New header cells: Wins / Losses / Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the existing header (A1) onto the new header cells
# so they share the same bold/border/centered style used by the rest of row 1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the team's win/loss/tie record on every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 63   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 99   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
